$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 65.666664
$ws.Range("I9").Value = 65.666664
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 65.666664
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 103.333336
$ws.Range("N9").ClearContents()
$ws.Range("H64").Value = 1500
$ws.Range("J64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("N64").Value = -1996
$ws.Range("H67").Value = 1500
$ws.Range("J67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("N67").Value = -3216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -84
$ws.Range("N4").ClearContents()
$ws.Range("H74").Value = 500
$ws.Range("I74").Value = 500
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 500
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 374
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 500
$ws.Range("I77").Value = 500
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1868
$ws.Range("N77").ClearContents()
$ws.Range("H132").Value = 1159.4
$ws.Range("I132").Value = 199.25
$ws.Range("K132").Value = 597.75
$ws.Range("M132").Value = 1932.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 814.125
$ws.Range("I16").Value = 760
$ws.Range("J16").Value = 904.3333
$ws.Range("K16").Value = 760
$ws.Range("L16").Value = 904.3333
$ws.Range("M16").Value = -473
$ws.Range("N16").Value = -1478.3333
$ws.Range("H31").Value = 11168.2
$ws.Range("I31").Value = 7699.8
$ws.Range("J31").Value = 12902.4
$ws.Range("K31").Value = 7699.8
$ws.Range("L31").Value = 12902.4
$ws.Range("M31").Value = -7404.8
$ws.Range("N31").Value = -13492.4
$ws.Range("H33").Value = 5981
$ws.Range("I33").Value = 979.8889
$ws.Range("K33").Value = 979.8889
$ws.Range("M33").Value = -600.8889
$ws.Range("H34").Value = 11168.2
$ws.Range("I34").Value = 7699.8
$ws.Range("J34").Value = 12902.4
$ws.Range("K34").Value = 7699.8
$ws.Range("L34").Value = 12902.4
$ws.Range("M34").Value = -7497.8
$ws.Range("N34").Value = -13306.4
$ws.Range("H36").Value = 5500
$ws.Range("I36").Value = 5500
$ws.Range("K36").Value = 5500
$ws.Range("M36").Value = -5112
$ws.Range("H40").Value = 5500
$ws.Range("I40").Value = 5500
$ws.Range("K40").Value = 5500
$ws.Range("M40").Value = -5340
$ws.Range("H58").Value = 1744.091
$ws.Range("I58").Value = 1642.7142
$ws.Range("J58").Value = 1921.5
$ws.Range("K58").Value = 1642.7142
$ws.Range("L58").Value = 1921.5
$ws.Range("M58").Value = -1439.7142
$ws.Range("N58").Value = -2327.5
$ws.Range("H113").Value = 814.125
$ws.Range("I113").Value = 760
$ws.Range("J113").Value = 904.3333
$ws.Range("K113").Value = 760
$ws.Range("L113").Value = 904.3333
$ws.Range("M113").Value = 1410
$ws.Range("N113").Value = -5244.3333
$ws.Range("H122").Value = 1162.5
$ws.Range("I122").Value = 1337
$ws.Range("J122").Value = 290
$ws.Range("K122").Value = 4011
$ws.Range("L122").Value = 870
$ws.Range("M122").Value = -1561
$ws.Range("N122").Value = -5770
$ws.Range("H134").Value = 1341.625
$ws.Range("I134").Value = 1322.1666
$ws.Range("K134").Value = 3966.4998
$ws.Range("M134").Value = -1431.4998
$ws.Range("H136").Value = 1744.091
$ws.Range("I136").Value = 1642.7142
$ws.Range("J136").Value = 1921.5
$ws.Range("K136").Value = 4928.142599999999
$ws.Range("L136").Value = 5764.5
$ws.Range("M136").Value = -2378.142599999999
$ws.Range("N136").Value = -10864.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 216.33333
$ws.Range("J103").Value = 212
$ws.Range("L103").Value = 636
$ws.Range("N103").Value = -2394
$ws.Range("H113").Value = 252.25
$ws.Range("J113").Value = 219.66667
$ws.Range("L113").Value = 659.00001
$ws.Range("N113").Value = -4999.00001
$ws.Range("H130").Value = 1900
$ws.Range("J130").Value = 1900
$ws.Range("L130").Value = 5700
$ws.Range("N130").Value = -15740
$ws.Range("H132").Value = 50
$ws.Range("I132").Value = 50
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 450
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 2080
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 10000000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 10000000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 10000000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -10000504
$ws.Range("H36").Value = 3464.2856
$ws.Range("I36").Value = 2636.3635
$ws.Range("J36").Value = 6500
$ws.Range("K36").Value = 2636.3635
$ws.Range("L36").Value = 6500
$ws.Range("M36").Value = -2151.3635
$ws.Range("N36").Value = -7470
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H113").Value = 4413
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 7107.6665
$ws.Range("J122").Value = 7161.5
$ws.Range("L122").Value = 21484.5
$ws.Range("N122").Value = -26384.5
$ws.Range("H126").Value = 12855.714
$ws.Range("I126").Value = 9997.5
$ws.Range("K126").Value = 29992.5
$ws.Range("M126").Value = -27522.5
$ws.Range("H132").Value = 761.3333
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 742
$ws.Range("K132").Value = 2400
$ws.Range("L132").Value = 2226
$ws.Range("M132").Value = 130
$ws.Range("N132").Value = -7286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 37750.223
$ws.Range("I40").Value = 24958
$ws.Range("K40").Value = 24958
$ws.Range("M40").Value = -24822
$ws.Range("H42").Value = 8342.667
$ws.Range("I42").Value = 5000
$ws.Range("K42").Value = 5000
$ws.Range("M42").Value = -4437
$ws.Range("H49").Value = 8342.667
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4853
$ws.Range("H54").Value = 30084
$ws.Range("J54").Value = 30084
$ws.Range("L54").Value = 30084
$ws.Range("N54").Value = -31372
$ws.Range("H61").Value = 1498.8
$ws.Range("J61").Value = 3500
$ws.Range("L61").Value = 3500
$ws.Range("N61").Value = -3904
$ws.Range("H113").Value = 1498.8
$ws.Range("J113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -7840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 26960
$ws.Range("I51").Value = 28509.334
$ws.Range("K51").Value = 28509.334
$ws.Range("M51").Value = -27999.334
$ws.Range("H96").Value = 15270.625
$ws.Range("I96").Value = 3033
$ws.Range("J96").Value = 35666.668
$ws.Range("K96").Value = 3033
$ws.Range("L96").Value = 35666.668
$ws.Range("M96").Value = -1660
$ws.Range("N96").Value = -38412.668
$ws.Range("H122").Value = 874.75
$ws.Range("I122").Value = 824.5
$ws.Range("J122").Value = 925
$ws.Range("K122").Value = 2473.5
$ws.Range("L122").Value = 2775
$ws.Range("M122").Value = -23.5
$ws.Range("N122").Value = -7675
$ws.Range("H125").Value = 64905
$ws.Range("J125").Value = 64905
$ws.Range("L125").Value = 64905
$ws.Range("N125").Value = -74745
$ws.Range("H133").Value = 120000
$ws.Range("J133").Value = 120000
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -130120
